$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188 (shifts existing rows 188:217 down to 189:218)
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with its values
$ws.Cells.Item(188, 1).Value = 10
$ws.Cells.Item(188, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(188, 3).Value = "La Araucanía"
$ws.Cells.Item(188, 4).Value = 45131
$ws.Cells.Item(188, 4).NumberFormat = $ws.Cells.Item(189, 4).NumberFormat
$ws.Cells.Item(188, 5).Value = 9
$ws.Cells.Item(188, 6).Value = 100112031
$ws.Cells.Item(188, 7).Value = "Poroto verde"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 55
$ws.Cells.Item(188, 11).Value = 25000
$ws.Cells.Item(188, 12).Value = 25000
$ws.Cells.Item(188, 13).Value = 25000
$ws.Cells.Item(188, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(188, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(188, 16).Value = 1000
$ws.Cells.Item(188, 17).Value = 25
$ws.Cells.Item(188, 18).Value = "Hortaliza"
